# Update "想去人数" (want-to-go count) numbers in column F across the
# workbook's four sheets, per the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1110
$ws1.Range("F3").Value  = 4797
$ws1.Range("F5").Value  = 1937
$ws1.Range("F6").Value  = 578
$ws1.Range("F7").Value  = 825
$ws1.Range("F8").Value  = 24
$ws1.Range("F9").Value  = 931
$ws1.Range("F10").Value = 1189
$ws1.Range("F11").Value = 1631
$ws1.Range("F12").Value = 856
$ws1.Range("F14").Value = 2061
$ws1.Range("F15").Value = 643
$ws1.Range("F17").Value = 541
$ws1.Range("F19").Value = 263
$ws1.Range("F20").Value = 132
$ws1.Range("F21").Value = 132
$ws1.Range("F23").Value = 1211
$ws1.Range("F25").Value = 2558
$ws1.Range("F27").Value = 15
$ws1.Range("F29").Value = 1672
$ws1.Range("F31").Value = 515
$ws1.Range("F34").Value = 4420
$ws1.Range("F35").Value = 72

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 4
$ws2.Range("F9").Value  = 50
$ws2.Range("F16").Value = 17
$ws2.Range("F18").Value = 149
$ws2.Range("F20").Value = 266
$ws2.Range("F24").Value = 1753
$ws2.Range("F26").Value = 212
$ws2.Range("F35").Value = 51

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1389
$ws3.Range("F5").Value = 1761
$ws3.Range("F7").Value = 440

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1389
$ws4.Range("F4").Value  = 1761
$ws4.Range("F6").Value  = 440
$ws4.Range("F7").Value  = 1110
$ws4.Range("F8").Value  = 4797
$ws4.Range("F9").Value  = 1937
$ws4.Range("F10").Value = 578
$ws4.Range("F13").Value = 24
$ws4.Range("F14").Value = 931
$ws4.Range("F15").Value = 1189
$ws4.Range("F16").Value = 1631
$ws4.Range("F18").Value = 50
$ws4.Range("F20").Value = 856
$ws4.Range("F22").Value = 2061
$ws4.Range("F23").Value = 643
$ws4.Range("F25").Value = 541
$ws4.Range("F27").Value = 263
$ws4.Range("F29").Value = 132
$ws4.Range("F30").Value = 132
$ws4.Range("F34").Value = 1211
$ws4.Range("F37").Value = 149
$ws4.Range("F38").Value = 2558
$ws4.Range("F39").Value = 266
$ws4.Range("F41").Value = 15
$ws4.Range("F43").Value = 1753
$ws4.Range("F44").Value = 1672
$ws4.Range("F45").Value = 515
$ws4.Range("F49").Value = 4420
$ws4.Range("F50").Value = 51
